$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.167.81'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.602.29'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.99'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.80%  '
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.14'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("D12").Value = '1.823.95'
$ws.Range("E12").Value = '  -0.50%  '
$ws.Range("D13").Value = '1.603.29'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("E14").Value = '  -1.14%  '
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '26.178.36'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.99'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  -0.90%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.96'
$ws.Range("E20").Value = '  +2.14%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.26'
$ws.Range("E22").Value = '  -2.51%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +11.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.55'
$ws.Range("E25").Value = '  +1.14%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("E27").Value = '  -7.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("E28").Value = '  -0.50%  '
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0489'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("E33").Value = '  -4.45%  '
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").Value = '1.141.95'
$ws.Range("E36").Value = '  +3.19%  '
$ws.Range("E37").Value = '  +6.45%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.786'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  -2.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.784'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.19'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = '1.738.31'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.21'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.50'
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.13'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  -11.66%  '
